$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handoff XLIFFs were regenerated -> status flips from "In Translation" to
# "Ready for handoff" and the generation timestamps move forward.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 02:46:09"

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 02:45:59"

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 02:46:09"

# The new status text is wider than the old one, so the status / language
# columns get re-sized to fit it.
$wsOverview.Range("E1:F1").ColumnWidth = 16.3
$wsZhCn.Range("C1").ColumnWidth = 16.3
$wsDeDe.Range("C1").ColumnWidth = 16.3
